$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.041.03"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").Value = "1.680.20"
$ws.Range("E3").Value = "  +0.87%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'215.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("E6").Value = "  -2.71%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  +2.11%  "

$ws.Range("D9").Value = "'21.42"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.78%  "

$ws.Range("E10").Value = "  +0.63%  "

$ws.Range("E11").Value = "  -0.81%  "

$ws.Range("D12").Value = "1.916.99"
$ws.Range("E12").Value = "  +0.83%  "

$ws.Range("D13").Value = "1.674.35"
$ws.Range("E13").Value = "  +0.47%  "

$ws.Range("E14").Value = "  +0.83%  "

$ws.Range("E15").Value = "  +1.64%  "

$ws.Range("D16").Value = "'66.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.33%  "

$ws.Range("D17").Value = "27.037.64"
$ws.Range("E17").Value = "  +0.44%  "

$ws.Range("D18").Value = "'8.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.71%  "

$ws.Range("D19").Value = "'235.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.49%  "

$ws.Range("E20").Value = "  +0.78%  "

$ws.Range("E21").Value = "  +0.01%  "

$ws.Range("D22").Value = "'4.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.24%  "

$ws.Range("D23").Value = "'9.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.93%  "

$ws.Range("E24").Value = "  -3.78%  "

$ws.Range("D25").Value = "'146.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.48%  "

$ws.Range("D26").Value = "'7.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.10%  "

$ws.Range("D27").Value = "'16.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.53%  "

$ws.Range("E28").Value = "  -2.26%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("E30").Value = "  +0.22%  "

$ws.Range("D32").Value = "'3.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.07%  "

$ws.Range("D33").Value = "1.539.95"
$ws.Range("E33").Value = "  +5.63%  "

$ws.Range("D34").Value = "'3.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.85%  "

$ws.Range("E35").Value = "  +5.64%  "

$ws.Range("E36").Value = "  -0.74%  "

$ws.Range("E37").Value = "  +1.35%  "

$ws.Range("D38").Value = "'0.915"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.20%  "

$ws.Range("E39").Value = "  +2.95%  "

$ws.Range("E40").Value = "  +6.47%  "

$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D42").Value = "'67.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.09%  "

$ws.Range("D43").Value = "'5.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.56%  "

$ws.Range("E44").Value = "  -0.53%  "

$ws.Range("D45").Value = "1.822.07"
$ws.Range("E45").Value = "  +0.79%  "

$ws.Range("D46").Value = "'0.779"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.38%  "

$ws.Range("D47").Value = "'90.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.09%  "

$ws.Range("E48").Value = "  +0.30%  "

$ws.Range("E49").Value = "  +2.15%  "

$ws.Range("D50").Value = "'8.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.42%  "

$ws.Range("E51").Value = "  +0.16%  "
